$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2086092715231788
$ws.Range("C2").Value = 0.5198675496688742
$ws.Range("J2").Value = 0.01324503311258278
$ws.Range("P2").Value = 0.1589403973509934
$ws.Range("S2").Value = 0.09933774834437085
$ws.Range("B3").Value = 0.006172839506172839
$ws.Range("C3").Value = 0.0308641975308642
$ws.Range("J3").Value = 0.04938271604938271
$ws.Range("P3").Value = 0.7160493827160493
$ws.Range("S3").Value = 0.1975308641975309
$ws.Range("J4").Value = 0.0392156862745098
$ws.Range("P4").Value = 0.6470588235294118
$ws.Range("S4").Value = 0.3137254901960784
$ws.Range("B6").Value = 0.07339449541284404
$ws.Range("D6").Value = 0.01376146788990826
$ws.Range("F6").Value = 0.07339449541284404
$ws.Range("J6").Value = 0.2522935779816514
$ws.Range("O6").Value = 0.004587155963302753
$ws.Range("Q6").Value = 0.1605504587155963
$ws.Range("R6").Value = 0.06880733944954129
$ws.Range("S6").Value = 0.3532110091743119
$ws.Range("B7").Value = 0.08666666666666667
$ws.Range("D7").Value = 0.04
$ws.Range("E7").Value = 0.006666666666666667
$ws.Range("F7").Value = 0.05333333333333334
$ws.Range("J7").Value = 0.1066666666666667
$ws.Range("O7").Value = 0.02666666666666667
$ws.Range("R7").Value = 0.08
$ws.Range("B8").Value = 0.100990099009901
$ws.Range("D8").Value = 0.0297029702970297
$ws.Range("F8").Value = 0.04752475247524753
$ws.Range("J8").Value = 0.1386138613861386
$ws.Range("O8").Value = 0.01584158415841584
$ws.Range("Q8").Value = 0.1782178217821782
$ws.Range("R8").Value = 0.09504950495049505
$ws.Range("S8").Value = 0.3940594059405941
$ws.Range("B9").Value = 0.1075949367088608
$ws.Range("D9").Value = 0.03164556962025317
$ws.Range("F9").Value = 0.0379746835443038
$ws.Range("J9").Value = 0.08227848101265822
$ws.Range("O9").Value = 0.01265822784810127
$ws.Range("Q9").Value = 0.2088607594936709
$ws.Range("R9").Value = 0.08860759493670886
$ws.Range("S9").Value = 0.4303797468354431
$ws.Range("B10").Value = 0.1006661732050333
$ws.Range("D10").Value = 0.01628423390081421
$ws.Range("F10").Value = 0.07549962990377498
$ws.Range("J10").Value = 0.1206513693560326
$ws.Range("O10").Value = 0.01332346410066617
$ws.Range("Q10").Value = 0.233160621761658
$ws.Range("R10").Value = 0.08734270910436713
$ws.Range("S10").Value = 0.3530717986676536
$ws.Range("G11").Value = 0.145748987854251
$ws.Range("J11").Value = 0.1174089068825911
$ws.Range("K11").Value = 0.1983805668016194
$ws.Range("L11").Value = 0.5182186234817814
$ws.Range("S11").Value = 0.02024291497975709
$ws.Range("G12").Value = 0.6716417910447762
$ws.Range("J12").Value = 0.2014925373134328
$ws.Range("K12").Value = 0.02985074626865672
$ws.Range("L12").Value = 0.05223880597014925
$ws.Range("S12").Value = 0.04477611940298507
$ws.Range("G13").Value = 0.5853658536585366
$ws.Range("J13").Value = 0.3414634146341464
$ws.Range("S13").Value = 0.07317073170731707
$ws.Range("F15").Value = 0.01244813278008299
$ws.Range("H15").Value = 0.2157676348547718
$ws.Range("I15").Value = 0.05394190871369295
$ws.Range("J15").Value = 0.3775933609958506
$ws.Range("K15").Value = 0.02904564315352697
$ws.Range("M15").Value = 0.01659751037344398
$ws.Range("O15").Value = 0.09958506224066389
$ws.Range("S15").Value = 0.1950207468879668
$ws.Range("F16").Value = 0.01058201058201058
$ws.Range("H16").Value = 0.1746031746031746
$ws.Range("I16").Value = 0.08994708994708994
$ws.Range("J16").Value = 0.3756613756613756
$ws.Range("K16").Value = 0.07936507936507936
$ws.Range("M16").Value = 0.02116402116402116
$ws.Range("O16").Value = 0.06878306878306878
$ws.Range("S16").Value = 0.1798941798941799
$ws.Range("F17").Value = 0.024
$ws.Range("H17").Value = 0.21
$ws.Range("I17").Value = 0.076
$ws.Range("J17").Value = 0.432
$ws.Range("K17").Value = 0.07199999999999999
$ws.Range("M17").Value = 0.018
$ws.Range("N17").Value = 0.002
$ws.Range("O17").Value = 0.066
$ws.Range("S17").Value = 0.1
$ws.Range("F18").Value = 0.01923076923076923
$ws.Range("H18").Value = 0.2019230769230769
$ws.Range("I18").Value = 0.08173076923076923
$ws.Range("J18").Value = 0.4519230769230769
$ws.Range("K18").Value = 0.08173076923076923
$ws.Range("M18").Value = 0.01923076923076923
$ws.Range("O18").Value = 0.05288461538461538
$ws.Range("S18").Value = 0.09134615384615384
$ws.Range("F19").Value = 0.008744038155802861
$ws.Range("H19").Value = 0.2225755166931637
$ws.Range("I19").Value = 0.05882352941176471
$ws.Range("J19").Value = 0.3934817170111288
$ws.Range("K19").Value = 0.09220985691573927
$ws.Range("N19").Value = 0.002384737678855326
$ws.Range("O19").Value = 0.0794912559618442
$ws.Range("S19").Value = 0.1240063593004769
